$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "GSM2463686-GSM2463688 "
$ws.Range("C7").Value = "N=17 (Severe`n PE) "
$ws.Range("C9").Value = "N=8 (Severe `nPE) "
$ws.Range("C15").Value = "N=6 PE`n Placental Tissue"
$ws.Range("C17").Value = "N=23 , PE placenta, `n(microarray)"
$ws.Range("C22").Value = "N=19, PE, placenta, `n(microarray) "
